$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 (item #9) - "الرصيد الحالي" (current balance) text value: 17:0 -> 16:0
$ws.Range("H12").Value = "16:0"

# "سعر البيع" (sale price) for item #9: 35 -> 70
$ws.Range("L12").Value = 70

# "عدد التعاملات" (number of transactions) for item #9: 1:0 -> 2:0
$ws.Range("N12").Value = "2:0"

# Total sale price row: 449.5 -> 484.5
$ws.Range("K13").Value = 484.5
